$d = $word.ActiveDocument

# --- Picture 1 (Section image) -> hyperlink ---
$url1 = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/GFA/GFA-42A-rooftop-communal-pavilion_section.jpg"
$shape1 = $d.InlineShapes.Item(1)
$range1 = $shape1.Range
$pos1 = $range1.Start
$shape1.Delete()
$insertAt1 = $d.Range($pos1, $pos1)
$d.Hyperlinks.Add($insertAt1, $url1, $null, $null, $url1) | Out-Null

# --- Picture 2 (Plan image) -> hyperlink ---
$url2 = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/GFA/GFA-42B-rooftop-communal-pavilion_plan.jpg?h=403&w=800"
$shape2 = $d.InlineShapes.Item(1)
$range2 = $shape2.Range
$pos2 = $range2.Start
$shape2.Delete()
$insertAt2 = $d.Range($pos2, $pos2)
$d.Hyperlinks.Add($insertAt2, $url2, $null, $null, $url2) | Out-Null

Write-Output "Replaced both pictures with hyperlinks."
